$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: dates for the next two days, formatted as "d-mmm"
$ws.Range("A1").Value = 44984
$ws.Range("B1").Value = 44985
$ws.Range("A1:B1").NumberFormat = "d-mmm"

# Fill column A first (rows 2-5), then column B (rows 2-5),
# so shared-string indices are allocated in that order.
$ws.Range("A2").Value = "Д. з. "
$ws.Range("A3").Value = "Информатика"
$ws.Range("A4").Value = "Структурка"
$ws.Range("A5").Value = "Эк. Геология"

$ws.Range("B2").Value = "Д. з."
$ws.Range("B3").Value = "Информатика"
$ws.Range("B4").Value = "ТФКП"
$ws.Range("B5").Value = "Структурка"

# Column widths (closest achievable values under this engine's
# pixel-quantized ColumnWidth implementation)
$ws.Columns.Item(1).ColumnWidth = 13.75
$ws.Columns.Item(2).ColumnWidth = 12.6

# Leave the active selection on B5, as in the edited workbook
$ws.Range("B5").Select() | Out-Null
